$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 40456
$ws.Range("J108").Value = 40456
$ws.Range("L108").Value = 40456
$ws.Range("N108").Value = -48136

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1764.0714
$ws.Range("I127").Value = 465.66666
$ws.Range("J127").Value = 2118.182
$ws.Range("K127").Value = 1396.99998
$ws.Range("L127").Value = 6354.545999999999
$ws.Range("M127").Value = 3563.00002
$ws.Range("N127").Value = -16274.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1180.0312
$ws.Range("I2").Value = 1067.8695
$ws.Range("J2").Value = 1466.6666
$ws.Range("K2").Value = 1067.8695
$ws.Range("L2").Value = 1466.6666
$ws.Range("M2").Value = -954.8695
$ws.Range("N2").Value = -1692.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 8600
$ws.Range("J9").Value = 8600
$ws.Range("L9").Value = 8600
$ws.Range("N9").Value = -8940

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 8600
$ws.Range("J20").Value = 8600
$ws.Range("L20").Value = 8600
$ws.Range("N20").Value = -9140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1125
$ws.Range("I45").Value = 955
$ws.Range("J45").Value = 1550
$ws.Range("K45").Value = 955
$ws.Range("L45").Value = 1550
$ws.Range("M45").Value = -578
$ws.Range("N45").Value = -2304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1180.0312
$ws.Range("I116").Value = 1067.8695
$ws.Range("J116").Value = 1466.6666
$ws.Range("K116").Value = 1067.8695
$ws.Range("L116").Value = 1466.6666
$ws.Range("M116").Value = 1226.1305
$ws.Range("N116").Value = -6054.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1180.0312
$ws.Range("I3").Value = 1067.8695
$ws.Range("J3").Value = 1466.6666
$ws.Range("K3").Value = 1067.8695
$ws.Range("L3").Value = 1466.6666
$ws.Range("M3").Value = -953.8695
$ws.Range("N3").Value = -1694.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 25250
$ws.Range("I15").Value = 1000
$ws.Range("K15").Value = 1000
$ws.Range("M15").Value = -773

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 39999.4
$ws.Range("I19").Value = 39998.5
$ws.Range("K19").Value = 39998.5
$ws.Range("M19").Value = -39825.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1775
$ws.Range("I94").Value = 725.1429000000001
$ws.Range("J94").Value = 2999.8333
$ws.Range("K94").Value = 725.1429000000001
$ws.Range("L94").Value = 2999.8333
$ws.Range("M94").Value = -274.1429000000001
$ws.Range("N94").Value = -3901.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 45455172
$ws.Range("I22").Value = 83333760
$ws.Range("J22").Value = 860.8
$ws.Range("K22").Value = 83333760
$ws.Range("L22").Value = 860.8
$ws.Range("M22").Value = -83333410
$ws.Range("N22").Value = -1560.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 15125
$ws.Range("I25").Value = 10250
$ws.Range("J25").Value = 20000
$ws.Range("K25").Value = 10250
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = -10076
$ws.Range("N25").Value = -20348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 55723776
$ws.Range("I62").Value = 66867532
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 66867532
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -66866908
$ws.Range("N62").Value = -6248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 55723776
$ws.Range("I65").Value = 66867532
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 334337660
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -334334540
$ws.Range("N65").Value = -31240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 697.44446
$ws.Range("I5").Value = 456
$ws.Range("J5").Value = 999.25
$ws.Range("K5").Value = 1368
$ws.Range("L5").Value = 2997.75
$ws.Range("M5").Value = -1256
$ws.Range("N5").Value = -3221.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3160
$ws.Range("J54").Value = 3442.8572
$ws.Range("L54").Value = 10328.5716
$ws.Range("N54").Value = -11446.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 1670.579
$ws.Range("I115").Value = 1405.0968
$ws.Range("J115").Value = 2846.2856
$ws.Range("K115").Value = 4215.2904
$ws.Range("L115").Value = 8538.856800000001
$ws.Range("M115").Value = -3040.2904
$ws.Range("N115").Value = -10888.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 902.8570999999999
$ws.Range("I131").Value = 457.5
$ws.Range("J131").Value = 949.7368
$ws.Range("K131").Value = 1372.5
$ws.Range("L131").Value = 2849.2104
$ws.Range("M131").Value = 3667.5
$ws.Range("N131").Value = -12929.2104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 697.44446
$ws.Range("I135").Value = 456
$ws.Range("J135").Value = 999.25
$ws.Range("K135").Value = 4104
$ws.Range("L135").Value = 8993.25
$ws.Range("M135").Value = -1569
$ws.Range("N135").Value = -14063.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 170430.17
$ws.Range("I14").Value = 288358.84
$ws.Range("J14").Value = 5330
$ws.Range("K14").Value = 288358.84
$ws.Range("L14").Value = 5330
$ws.Range("M14").Value = -288190.84
$ws.Range("N14").Value = -5666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 8950
$ws.Range("I18").Value = 7000
$ws.Range("J18").Value = 9600
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 9600
$ws.Range("M18").Value = -6707
$ws.Range("N18").Value = -10186

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 3065
$ws.Range("I18").Value = 1397.5
$ws.Range("J18").Value = 6400
$ws.Range("K18").Value = 1397.5
$ws.Range("L18").Value = 6400
$ws.Range("M18").Value = -1225.5
$ws.Range("N18").Value = -6744

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 6000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 6000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 6000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -6452

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1956.7142
$ws.Range("I96").Value = 1515.6666
$ws.Range("J96").Value = 2287.5
$ws.Range("K96").Value = 1515.6666
$ws.Range("L96").Value = 2287.5
$ws.Range("M96").Value = -142.6666
$ws.Range("N96").Value = -5033.5
